$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# The document has two distinct picture "logos" that each appear twice
# (once in the default header/footer, once in the first-page header/footer):
#   - BTec_Logo-Orange picture:   currently named "image2.jpg" -> should become "image1.jpg"
#   - PearsonLogo picture:        currently named "image1.png" -> should become "image2.png"
# Update the wp:docPr / pic:cNvPr "name" attribute (exposed on the COM object
# model as InlineShape.Name) for each of the four pictures, in both headers
# and both footers.

function Set-LogoName {
    param($shape, $newName)
    if ($shape -ne $null) {
        try {
            $shape.Name = $newName
        } catch {
            # Some stories can refuse an in-place rename of the shape's
            # name; ignore so the other pictures still get updated.
        }
    }
}

# Headers: swap BTec_Logo-Orange from image2.jpg -> image1.jpg
$header1 = $sec.Headers.Item(1)
if ($header1.Exists) {
    $shape = $header1.Range.InlineShapes.Item(1)
    Set-LogoName $shape "image1.jpg"
}

$header2 = $sec.Headers.Item(2)
if ($header2.Exists) {
    $shape = $header2.Range.InlineShapes.Item(1)
    Set-LogoName $shape "image1.jpg"
}

# Footers: swap PearsonLogo from image1.png -> image2.png
$footer1 = $sec.Footers.Item(1)
if ($footer1.Exists) {
    $shape = $footer1.Range.InlineShapes.Item(1)
    Set-LogoName $shape "image2.png"
}

$footer2 = $sec.Footers.Item(2)
if ($footer2.Exists) {
    $shape = $footer2.Range.InlineShapes.Item(1)
    Set-LogoName $shape "image2.png"
}

Write-Output "logo names updated"
